$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Q_toBoil value (B3) with the new uncertainty-adjusted figure.
# B4 (T_boil) holds formula =B3/B2 and will recalculate automatically.
$ws.Range("B3").Value = 511054.76033279998

$excel.Calculate()
